$d = $word.ActiveDocument

# 1. Make the two empty runs (before "DRY" body text and before the
#    "What is SOLID Design Principles" heading) bold, matching the
#    surrounding heading style.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Bold = 1
$p2.Range.Font.BoldBi = 1

$p4 = $d.Paragraphs.Item(4)
$p4.Range.Font.Bold = 1
$p4.Range.Font.BoldBi = 1

# 2. Merge "What is SOLID " + "Design Principles" into a single run of
#    text (no functional text change, just collapses the run split).
$d.Content.Find.Execute("What is SOLID Design Principles", $true, $false, $false, $false, $false, $true, 1, $false, "What is SOLID Design Principles", 2)

# 3. Swap the two "Easier to ..." bullet points.
$d.Content.Find.Execute("Easier to understand", $true, $false, $false, $false, $false, $true, 1, $false, "Easier to TEMP_PLACEHOLDER", 2)
$d.Content.Find.Execute("Easier to maintain", $true, $false, $false, $false, $false, $true, 1, $false, "Easier to understand", 2)
$d.Content.Find.Execute("Easier to TEMP_PLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "Easier to maintain", 2)

# 4. Add spaces around the slash in the OCP heading.
$d.Content.Find.Execute("2. Open/Closed Principle (OCP)", $true, $false, $false, $false, $false, $true, 1, $false, "2. Open / Closed Principle (OCP)", 2)
